$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 576.9091
$ws.Range("I8").Value = 44.5
$ws.Range("J8").Value = 1996.6666
$ws.Range("K8").Value = 133.5
$ws.Range("L8").Value = 5989.9998
$ws.Range("M8").Value = 5.5
$ws.Range("N8").Value = -6267.9998
$ws.Range("H100").Value = 93494.17999999999
$ws.Range("I100").Value = 101843.6
$ws.Range("K100").Value = 101843.6
$ws.Range("M100").Value = -101302.6
$ws.Range("H113").Value = 3578.6667
$ws.Range("I113").Value = 2473.3333
$ws.Range("K113").Value = 2473.3333
$ws.Range("M113").Value = 780.6667000000002
$ws.Range("H116").Value = 7817.28
$ws.Range("I116").Value = 6918.3687
$ws.Range("J116").Value = 10663.833
$ws.Range("K116").Value = 6918.3687
$ws.Range("L116").Value = 10663.833
$ws.Range("M116").Value = -3476.3687
$ws.Range("N116").Value = -17547.833
$ws.Range("H118").Value = 555.2143
$ws.Range("I118").Value = 580.2308
$ws.Range("J118").Value = 230
$ws.Range("K118").Value = 1740.6924
$ws.Range("L118").Value = 690
$ws.Range("M118").Value = -83.69240000000013
$ws.Range("N118").Value = -4004
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("H132").Value = 2660.6875
$ws.Range("I132").Value = 2708.4
$ws.Range("J132").Value = 1945
$ws.Range("K132").Value = 8125.200000000001
$ws.Range("L132").Value = 5835
$ws.Range("M132").Value = -5595.200000000001
$ws.Range("N132").Value = -10895
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3632.82
$ws.Range("I2").Value = 3578.8
$ws.Range("J2").Value = 3758.8667
$ws.Range("K2").Value = 3578.8
$ws.Range("L2").Value = 3758.8667
$ws.Range("M2").Value = -3465.8
$ws.Range("N2").Value = -3984.8667
$ws.Range("H32").Value = 7678.278
$ws.Range("I32").Value = 7052.125
$ws.Range("J32").Value = 12687.5
$ws.Range("K32").Value = 7052.125
$ws.Range("L32").Value = 12687.5
$ws.Range("M32").Value = -6765.125
$ws.Range("N32").Value = -13261.5
$ws.Range("H43").Value = 10265.333
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10265.333
$ws.Range("K43").Value = 0
$ws.Range("N43").Value = -10891.333
$ws.Range("H97").Value = 1185.9412
$ws.Range("I97").Value = 1060.2963
$ws.Range("J97").Value = 1670.5714
$ws.Range("K97").Value = 1060.2963
$ws.Range("L97").Value = 1670.5714
$ws.Range("M97").Value = -564.2963
$ws.Range("N97").Value = -2662.5714
$ws.Range("H116").Value = 3632.82
$ws.Range("I116").Value = 3578.8
$ws.Range("J116").Value = 3758.8667
$ws.Range("K116").Value = 3578.8
$ws.Range("L116").Value = 3758.8667
$ws.Range("M116").Value = -1284.8
$ws.Range("N116").Value = -8346.8667
$ws.Range("H132").Value = 934.87177
$ws.Range("I132").Value = 912.25
$ws.Range("J132").Value = 1038.2858
$ws.Range("K132").Value = 2736.75
$ws.Range("L132").Value = 3114.8574
$ws.Range("M132").Value = -206.75
$ws.Range("N132").Value = -8174.857400000001
$ws.Range("L43").Value = 10265.333
$ws.Range("M43").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3632.82
$ws.Range("I3").Value = 3578.8
$ws.Range("J3").Value = 3758.8667
$ws.Range("K3").Value = 3578.8
$ws.Range("L3").Value = 3758.8667
$ws.Range("M3").Value = -3464.8
$ws.Range("N3").Value = -3986.8667
$ws.Range("H94").Value = 1968.0526
$ws.Range("I94").Value = 1513.8572
$ws.Range("J94").Value = 3239.8
$ws.Range("K94").Value = 1513.8572
$ws.Range("L94").Value = 3239.8
$ws.Range("M94").Value = -1062.8572
$ws.Range("N94").Value = -4141.8
$ws.Range("H99").Value = 5029.525
$ws.Range("I99").Value = 5023.838
$ws.Range("J99").Value = 5099.6665
$ws.Range("K99").Value = 5023.838
$ws.Range("L99").Value = 5099.6665
$ws.Range("M99").Value = -3525.838
$ws.Range("N99").Value = -8095.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2972.9722
$ws.Range("I16").Value = 3073.5806
$ws.Range("J16").Value = 2349.2
$ws.Range("K16").Value = 3073.5806
$ws.Range("L16").Value = 2349.2
$ws.Range("M16").Value = -2786.5806
$ws.Range("N16").Value = -2923.2
$ws.Range("H23").Value = 3227.2727
$ws.Range("J23").Value = 5500
$ws.Range("L23").Value = 5500
$ws.Range("N23").Value = -5980
$ws.Range("H27").Value = 3227.2727
$ws.Range("J27").Value = 5500
$ws.Range("L27").Value = 5500
$ws.Range("N27").Value = -5884
$ws.Range("H31").Value = 4379.846
$ws.Range("I31").Value = 3056.6667
$ws.Range("J31").Value = 4776.8
$ws.Range("K31").Value = 3056.6667
$ws.Range("L31").Value = 4776.8
$ws.Range("M31").Value = -2761.6667
$ws.Range("N31").Value = -5366.8
$ws.Range("H34").Value = 4379.846
$ws.Range("I34").Value = 3056.6667
$ws.Range("J34").Value = 4776.8
$ws.Range("K34").Value = 3056.6667
$ws.Range("L34").Value = 4776.8
$ws.Range("M34").Value = -2854.6667
$ws.Range("N34").Value = -5180.8
$ws.Range("H35").Value = 846.5
$ws.Range("I35").Value = 443
$ws.Range("K35").Value = 443
$ws.Range("M35").Value = -149
$ws.Range("H62").Value = 68173.375
$ws.Range("I62").Value = 145686.58
$ws.Range("J62").Value = 7885.3335
$ws.Range("K62").Value = 145686.58
$ws.Range("L62").Value = 7885.3335
$ws.Range("M62").Value = -145062.58
$ws.Range("N62").Value = -9133.333500000001
$ws.Range("H65").Value = 68173.375
$ws.Range("I65").Value = 145686.58
$ws.Range("J65").Value = 7885.3335
$ws.Range("K65").Value = 728432.8999999999
$ws.Range("L65").Value = 39426.6675
$ws.Range("M65").Value = -725312.8999999999
$ws.Range("N65").Value = -45666.6675
$ws.Range("H107").Value = 949.1539
$ws.Range("I107").Value = 953.6111
$ws.Range("J107").Value = 939.125
$ws.Range("K107").Value = 953.6111
$ws.Range("L107").Value = 939.125
$ws.Range("M107").Value = 966.3889
$ws.Range("N107").Value = -4779.125
$ws.Range("H113").Value = 2972.9722
$ws.Range("I113").Value = 3073.5806
$ws.Range("J113").Value = 2349.2
$ws.Range("K113").Value = 3073.5806
$ws.Range("L113").Value = 2349.2
$ws.Range("M113").Value = -903.5805999999998
$ws.Range("N113").Value = -6689.2
$ws.Range("H132").Value = 28242.06
$ws.Range("I132").Value = 17281.037
$ws.Range("J132").Value = 41109.348
$ws.Range("K132").Value = 51843.111
$ws.Range("L132").Value = 123328.044
$ws.Range("M132").Value = -49313.111
$ws.Range("N132").Value = -128388.044
$ws.Range("H134").Value = 5120.8647
$ws.Range("I134").Value = 3949.9524
$ws.Range("J134").Value = 6657.6875
$ws.Range("K134").Value = 11849.8572
$ws.Range("L134").Value = 19973.0625
$ws.Range("M134").Value = -9314.8572
$ws.Range("N134").Value = -25043.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1022.4
$ws.Range("J29").Value = 2514.5
$ws.Range("L29").Value = 7543.5
$ws.Range("N29").Value = -8097.5
$ws.Range("H131").Value = 5970
$ws.Range("I131").Value = 1655.3636
$ws.Range("J131").Value = 8033.522
$ws.Range("K131").Value = 4966.0908
$ws.Range("L131").Value = 24100.566
$ws.Range("M131").Value = 73.90920000000006
$ws.Range("N131").Value = -34180.566
$ws.Range("H132").Value = 2859350.2
$ws.Range("I132").Value = 1535
$ws.Range("J132").Value = 10003888
$ws.Range("K132").Value = 13815
$ws.Range("L132").Value = 90034992
$ws.Range("M132").Value = -11285
$ws.Range("N132").Value = -90040052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66.95
$ws.Range("I2").Value = 24.307692
$ws.Range("K2").Value = 24.307692
$ws.Range("M2").Value = 88.692308
$ws.Range("H80").Value = 2828.8333
$ws.Range("I80").Value = 2814.6
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 2814.6
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -1816.6
$ws.Range("N80").Value = -4896
$ws.Range("H83").Value = 2828.8333
$ws.Range("I83").Value = 2814.6
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 14073
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -9081
$ws.Range("N83").Value = -24484
$ws.Range("H132").Value = 10046.692
$ws.Range("I132").Value = 11234.143
$ws.Range("J132").Value = 5059.4
$ws.Range("K132").Value = 33702.429
$ws.Range("L132").Value = 15178.2
$ws.Range("M132").Value = -31172.429
$ws.Range("N132").Value = -20238.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3222.1738
$ws.Range("I68").Value = 2681.1667
$ws.Range("K68").Value = 2681.1667
$ws.Range("M68").Value = -1932.1667
$ws.Range("H71").Value = 3222.1738
$ws.Range("I71").Value = 2681.1667
$ws.Range("K71").Value = 13405.8335
$ws.Range("M71").Value = -9661.833500000001
$ws.Range("H132").Value = 7037.4736
$ws.Range("I132").Value = 6781.533
$ws.Range("J132").Value = 7997.25
$ws.Range("K132").Value = 20344.599
$ws.Range("L132").Value = 23991.75
$ws.Range("M132").Value = -17814.599
$ws.Range("N132").Value = -29051.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 4450.375
$ws.Range("I17").Value = 4450.375
$ws.Range("K17").Value = 4450.375
$ws.Range("M17").Value = -4278.375
$ws.Range("H40").Value = 34171
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 4839.5
$ws.Range("I62").Value = 3564.8333
$ws.Range("J62").Value = 6751.5
$ws.Range("K62").Value = 3564.8333
$ws.Range("L62").Value = 6751.5
$ws.Range("M62").Value = -2940.8333
$ws.Range("N62").Value = -7999.5
$ws.Range("H65").Value = 4839.5
$ws.Range("I65").Value = 3564.8333
$ws.Range("J65").Value = 6751.5
$ws.Range("K65").Value = 17824.1665
$ws.Range("L65").Value = 33757.5
$ws.Range("M65").Value = -14704.1665
$ws.Range("N65").Value = -39997.5
$ws.Range("H107").Value = 1620.4572
$ws.Range("I107").Value = 1419.8889
$ws.Range("J107").Value = 2297.375
$ws.Range("K107").Value = 4259.6667
$ws.Range("L107").Value = 6892.125
$ws.Range("M107").Value = -2339.6667
$ws.Range("N107").Value = -10732.125
$ws.Range("H122").Value = 1475.3214
$ws.Range("I122").Value = 713
$ws.Range("J122").Value = 6049.25
$ws.Range("K122").Value = 2139
$ws.Range("L122").Value = 18147.75
$ws.Range("M122").Value = 311
$ws.Range("N122").Value = -23047.75
$ws.Range("H136").Value = 1651.8776
$ws.Range("I136").Value = 1096.1666
$ws.Range("J136").Value = 3190.7693
$ws.Range("K136").Value = 3288.4998
$ws.Range("L136").Value = 9572.3079
$ws.Range("M136").Value = -738.4998000000001
$ws.Range("N136").Value = -14672.3079
